$wb = $excel.ActiveWorkbook

# This script rewrites the static "current market price" snapshot columns
# (H:N) that a scheduled data-refresh run recomputed for a handful of leves
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets. No formulas are involved;
# every value below is a hard-coded number written by the external price feed.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 106.13333
$ws.Range("I2").Value = 106.42857
$ws.Range("K2").Value = 106.42857
$ws.Range("M2").Value = 6.571430000000007
$ws.Range("H19").Value = 666.96
$ws.Range("I19").Value = 894
$ws.Range("J19").Value = 515.6
$ws.Range("K19").Value = 894
$ws.Range("L19").Value = 515.6
$ws.Range("M19").Value = -719
$ws.Range("N19").Value = -865.6
$ws.Range("H33").Value = 122.818184
$ws.Range("I33").Value = 75.2
$ws.Range("K33").Value = 75.2
$ws.Range("M33").Value = 153.8
$ws.Range("H40").Value = 6063.9
$ws.Range("I40").Value = 5074.769
$ws.Range("J40").Value = 6820.294
$ws.Range("K40").Value = 5074.769
$ws.Range("L40").Value = 6820.294
$ws.Range("M40").Value = -4899.769
$ws.Range("N40").Value = -7170.294
$ws.Range("H51").Value = 9999
$ws.Range("J51").Value = 9999
$ws.Range("L51").Value = 9999
$ws.Range("N51").Value = -10967
$ws.Range("H138").Value = 3673.5
$ws.Range("J138").Value = 4748.5
$ws.Range("L138").Value = 14245.5
$ws.Range("N138").Value = -24525.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3490.1538
$ws.Range("I122").Value = 3490.1538
$ws.Range("K122").Value = 10470.4614
$ws.Range("M122").Value = -8020.4614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1680.6364
$ws.Range("I20").Value = 1176.3334
$ws.Range("J20").Value = 2285.8
$ws.Range("K20").Value = 1176.3334
$ws.Range("L20").Value = 2285.8
$ws.Range("M20").Value = -929.3334
$ws.Range("N20").Value = -2779.8
$ws.Range("H29").Value = 800
$ws.Range("I29").Value = 800
$ws.Range("K29").Value = 800
$ws.Range("M29").Value = -511
$ws.Range("H36").Value = 3500
$ws.Range("I36").Value = 166.66667
$ws.Range("K36").Value = 166.66667
$ws.Range("M36").Value = 367.33333
$ws.Range("H75").Value = 9124
$ws.Range("I75").Value = 9124
$ws.Range("K75").Value = 9124
$ws.Range("M75").Value = -8188
$ws.Range("H78").Value = 9124
$ws.Range("I78").Value = 9124
$ws.Range("K78").Value = 27372
$ws.Range("M78").Value = -22692
$ws.Range("H105").Value = 1509.6
$ws.Range("I105").Value = 1494.7778
$ws.Range("J105").Value = 1531.8334
$ws.Range("K105").Value = 1494.7778
$ws.Range("L105").Value = 1531.8334
$ws.Range("M105").Value = 252.2221999999999
$ws.Range("N105").Value = -5025.8334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 1570
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 1950
$ws.Range("K62").Value = 1000
$ws.Range("L62").Value = 1950
$ws.Range("M62").Value = -376
$ws.Range("N62").Value = -3198
$ws.Range("H65").Value = 1570
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 1950
$ws.Range("K65").Value = 5000
$ws.Range("L65").Value = 9750
$ws.Range("M65").Value = -1880
$ws.Range("N65").Value = -15990
$ws.Range("H70").Value = 45000
$ws.Range("J70").Value = 45000
$ws.Range("L70").Value = 45000
$ws.Range("N70").Value = -45630
$ws.Range("H73").Value = 45000
$ws.Range("J73").Value = 45000
$ws.Range("L73").Value = 45000
$ws.Range("N73").Value = -47184
$ws.Range("H86").Value = 7591.2
$ws.Range("I86").Value = 7414.5713
$ws.Range("J86").Value = 8003.3335
$ws.Range("K86").Value = 7414.5713
$ws.Range("L86").Value = 8003.3335
$ws.Range("M86").Value = -6291.5713
$ws.Range("N86").Value = -10249.3335
$ws.Range("H89").Value = 7591.2
$ws.Range("I89").Value = 7414.5713
$ws.Range("J89").Value = 8003.3335
$ws.Range("K89").Value = 37072.85649999999
$ws.Range("L89").Value = 40016.6675
$ws.Range("M89").Value = -31456.85649999999
$ws.Range("N89").Value = -51248.6675
$ws.Range("H99").Value = 1867.1111
$ws.Range("I99").Value = 1884
$ws.Range("K99").Value = 1884
$ws.Range("M99").Value = -386
$ws.Range("H126").Value = 1867.1111
$ws.Range("I126").Value = 1884
$ws.Range("K126").Value = 5652
$ws.Range("M126").Value = -3182

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 83.111115
$ws.Range("I17").Value = 24
$ws.Range("J17").Value = 130.4
$ws.Range("K17").Value = 72
$ws.Range("L17").Value = 391.2
$ws.Range("M17").Value = 97
$ws.Range("N17").Value = -729.2
$ws.Range("H68").Value = 3599.2
$ws.Range("J68").Value = 4667.6665
$ws.Range("L68").Value = 14002.9995
$ws.Range("N68").Value = -15624.9995
$ws.Range("H71").Value = 3599.2
$ws.Range("J71").Value = 4667.6665
$ws.Range("L71").Value = 42008.9985
$ws.Range("N71").Value = -50120.9985
$ws.Range("H80").Value = 4262.7085
$ws.Range("J80").Value = 4757.857
$ws.Range("L80").Value = 14273.571
$ws.Range("N80").Value = -16145.571
$ws.Range("H83").Value = 4262.7085
$ws.Range("J83").Value = 4757.857
$ws.Range("L83").Value = 42820.713
$ws.Range("N83").Value = -52180.713
$ws.Range("H132").Value = 4444
$ws.Range("I132").Value = 4444
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 39996
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -37466
$ws.Range("N132").Value = $null
$ws.Range("H140").Value = 2435.5715
$ws.Range("I140").Value = 2008.1666
$ws.Range("K140").Value = 6024.4998
$ws.Range("M140").Value = -844.4997999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8998
$ws.Range("I5").Value = 8998
$ws.Range("K5").Value = 8998
$ws.Range("M5").Value = -8886
$ws.Range("H97").Value = 437.33334
$ws.Range("I97").Value = 437.33334
$ws.Range("K97").Value = 437.33334
$ws.Range("M97").Value = 58.66665999999998
$ws.Range("H126").Value = 1750
$ws.Range("I126").Value = 1750
$ws.Range("K126").Value = 5250
$ws.Range("M126").Value = -2780

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 782.5714
$ws.Range("J16").Value = 486
$ws.Range("L16").Value = 486
$ws.Range("N16").Value = -826
$ws.Range("H46").Value = 4353.5
$ws.Range("I46").Value = 3631.4
$ws.Range("K46").Value = 3631.4
$ws.Range("M46").Value = -3443.4
$ws.Range("H122").Value = 4600.4
$ws.Range("I122").Value = 4600.4
$ws.Range("K122").Value = 13801.2
$ws.Range("M122").Value = -11351.2
$ws.Range("H132").Value = 9566.091
$ws.Range("I132").Value = 5856
$ws.Range("K132").Value = 17568
$ws.Range("M132").Value = -15038

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7250.25
$ws.Range("I81").Value = 7250.25
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 14500.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -13439.5
$ws.Range("N81").Value = $null
$ws.Range("H84").Value = 7250.25
$ws.Range("I84").Value = 7250.25
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 72502.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -67198.5
$ws.Range("N84").Value = $null
